$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B-E on this sheet store plain text (inline strings) even when the text looks
# like a number (e.g. "215.01" is the literal price text, not a numeric value). Briefly
# switch affected "D" cells to a text format while assigning so Excel keeps the value as
# a string instead of silently parsing it into a real number, then restore "General".

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '26.959.53'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value2 = '  -0.28%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '1.676.73'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value2 = '  +0.10%  '

$ws.Range("E4").Value2 = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '215.01'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value2 = '  -0.71%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '0.518'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value2 = '  -2.50%  '

$ws.Range("E8").Value2 = '  -1.14%  '

$ws.Range("E9").Value2 = '  -0.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '21.10'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value2 = '  +4.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '0.0888'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value2 = '  -0.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '1.911.89'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value2 = '  +0.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '1.674.75'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value2 = '  -0.67%  '

$ws.Range("E14").Value2 = '  +0.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '0.531'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value2 = '  +1.87%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '65.77'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value2 = '  -0.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '8.23'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value2 = '  +5.81%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '26.965.94'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value2 = '  -0.30%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '236.59'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value2 = '  +1.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '0.0₃0734'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value2 = '  -0.43%  '

$ws.Range("E21").Value2 = '  -0.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '4.44'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value2 = '  -0.57%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '9.20'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value2 = '  -0.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '2.13'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value2 = '  -4.42%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '146.88'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value2 = '  +0.77%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '7.23'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value2 = '  +1.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '16.06'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value2 = '  +1.09%  '

$ws.Range("E28").Value2 = '  -2.88%  '

$ws.Range("E30").Value2 = '  +0.54%  '

$ws.Range("E31").Value2 = '  -1.42%  '

$ws.Range("E32").Value2 = '  +0.25%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '1.488.89'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value2 = '  +2.49%  '

$ws.Range("E34").Value2 = '  +0.68%  '

$ws.Range("E35").Value2 = '  +4.46%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '2.41'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value2 = '  +0.12%  '

$ws.Range("E37").Value2 = '  +3.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '0.0175'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value2 = '  +3.50%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '0.911'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value2 = '  +1.16%  '

$ws.Range("B40").Value2 = 'WEMIXToken'
$ws.Range("C40").Value2 = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '1.03'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value2 = '  +4.92%  '

$ws.Range("B41").Value2 = 'FraxShare'
$ws.Range("C41").Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '5.76'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value2 = '  -5.00%  '

$ws.Range("E42").Value2 = '  +0.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '67.46'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value2 = '  +2.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '2.28'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value2 = '  -1.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '1.818.06'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value2 = '  +0.33%  '

$ws.Range("E46").Value2 = '  -0.52%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '90.51'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value2 = '  -0.15%  '

$ws.Range("E48").Value2 = '  -0.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '0.104'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value2 = '  +3.06%  '

$ws.Range("E50").Value2 = '  +0.23%  '

$ws.Range("E51").Value2 = '  +1.20%  '
